$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume values from the latest scrape.
# NumberFormat is temporarily forced to text ("@") before assigning the
# value so numeric-looking strings (e.g. "361.07") are stored as text,
# matching the original inline-string cell type; ClearFormats() then
# removes the temporary format so no visible style is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '52.161.89'
Set-TextValue $ws.Range('E2') '  -0.16%  '
Set-TextValue $ws.Range('D3') '2.845.90'
Set-TextValue $ws.Range('E3') '  +1.90%  '
Set-TextValue $ws.Range('E4') '  +0.02%  '
Set-TextValue $ws.Range('D5') '361.07'
Set-TextValue $ws.Range('E5') '  +5.79%  '
Set-TextValue $ws.Range('D6') '113.52'
Set-TextValue $ws.Range('E6') '  -3.11%  '
Set-TextValue $ws.Range('D7') '0.576'
Set-TextValue $ws.Range('E7') '  +4.19%  '
Set-TextValue $ws.Range('E8') '  +0.00%  '
Set-TextValue $ws.Range('D9') '0.601'
Set-TextValue $ws.Range('E9') '  +3.72%  '
Set-TextValue $ws.Range('D10') '41.55'
Set-TextValue $ws.Range('E10') '  -1.53%  '
Set-TextValue $ws.Range('D11') '0.0863'
Set-TextValue $ws.Range('E11') '  -0.89%  '
Set-TextValue $ws.Range('E12') '  +1.11%  '
Set-TextValue $ws.Range('D13') '20.00'
Set-TextValue $ws.Range('E13') '  -0.60%  '
Set-TextValue $ws.Range('E14') '  +1.94%  '
Set-TextValue $ws.Range('D15') '3.294.94'
Set-TextValue $ws.Range('E15') '  +1.84%  '
Set-TextValue $ws.Range('D16') '2.848.37'
Set-TextValue $ws.Range('E16') '  +2.52%  '
Set-TextValue $ws.Range('D17') '0.904'
Set-TextValue $ws.Range('E17') '  +1.89%  '
Set-TextValue $ws.Range('D18') '51.957.51'
Set-TextValue $ws.Range('E18') '  -0.23%  '
Set-TextValue $ws.Range('E19') '  +8.95%  '
Set-TextValue $ws.Range('D20') '3.17'
Set-TextValue $ws.Range('E20') '  -2.47%  '
Set-TextValue $ws.Range('D21') '13.59'
Set-TextValue $ws.Range('E21') '  +1.60%  '
Set-TextValue $ws.Range('E22') '  +0.51%  '
Set-TextValue $ws.Range('D23') '70.33'
Set-TextValue $ws.Range('E23') '  -0.04%  '
Set-TextValue $ws.Range('D24') '267.94'
Set-TextValue $ws.Range('E24') '  -3.79%  '
Set-TextValue $ws.Range('E25') '  +0.32%  '
Set-TextValue $ws.Range('D26') '27.21'
Set-TextValue $ws.Range('E26') '  +1.02%  '
Set-TextValue $ws.Range('D27') '0.999'
Set-TextValue $ws.Range('E27') '  +0.04%  '
Set-TextValue $ws.Range('E28') '  +2.00%  '
Set-TextValue $ws.Range('E29') '  +1.42%  '
Set-TextValue $ws.Range('D30') '53.64'
Set-TextValue $ws.Range('E30') '  +6.36%  '
Set-TextValue $ws.Range('E31') '  -1.00%  '
Set-TextValue $ws.Range('D32') '0.0460'
Set-TextValue $ws.Range('E32') '  +23.66%  '
Set-TextValue $ws.Range('D33') '34.20'
Set-TextValue $ws.Range('E33') '  -2.30%  '
Set-TextValue $ws.Range('D34') '5.90'
Set-TextValue $ws.Range('E34') '  +3.18%  '
Set-TextValue $ws.Range('D35') '5.42'
Set-TextValue $ws.Range('E35') '  +8.67%  '
Set-TextValue $ws.Range('D36') '0.0845'
Set-TextValue $ws.Range('E36') '  +2.21%  '
Set-TextValue $ws.Range('E37') '  -0.01%  '
Set-TextValue $ws.Range('E38') '  +0.23%  '
Set-TextValue $ws.Range('E39') '  -2.61%  '
Set-TextValue $ws.Range('D40') '18.35'
Set-TextValue $ws.Range('E40') '  -3.60%  '
Set-TextValue $ws.Range('D41') '23.77'
Set-TextValue $ws.Range('E41') '  +1.43%  '
Set-TextValue $ws.Range('D42') '0.118'
Set-TextValue $ws.Range('E42') '  +1.38%  '
Set-TextValue $ws.Range('D43') '128.50'
Set-TextValue $ws.Range('E43') '  +1.07%  '
Set-TextValue $ws.Range('E44') '  -7.67%  '
Set-TextValue $ws.Range('E45') '  -3.03%  '
Set-TextValue $ws.Range('D46') '2.114.35'
Set-TextValue $ws.Range('E46') '  +0.33%  '
Set-TextValue $ws.Range('E47') '  +1.28%  '
Set-TextValue $ws.Range('E49') '  +8.73%  '
Set-TextValue $ws.Range('E50') '  +5.38%  '
Set-TextValue $ws.Range('D51') '9.06'
Set-TextValue $ws.Range('E51') '  +1.10%  '
